$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '20.597.01'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.47%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.468.69'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.95%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9621'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.85%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '281.49'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3727'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("E7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3209'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.73%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.83'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.10%  '
$ws.Range("E9").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.25%  '
$ws.Range("E10").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("E11").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.37%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.653'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.41%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.40'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.42%  '
$ws.Range("E14").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.475.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.15%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001042'
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05788'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("E18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.96'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.72%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9643'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.61%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.720'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.90'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.42%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.32'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("E23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.308'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.97%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '20.648.59'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.334'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '138.05'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("E27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.638.28'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.97%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.04'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.10%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.975'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.390'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.90%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8452'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.64%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.639'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +25.92%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07873'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("E35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06102'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.99%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.956'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.76'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.30%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02087'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.94%  '
$ws.Range("E39").ClearFormats()
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.137'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.29%  '
$ws.Range("E40").ClearFormats()
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9726'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.75%  '
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1910'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.622'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -9.63%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5469'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.18%  '
$ws.Range("E44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.62'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.41%  '
$ws.Range("E45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.600'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.45%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '121.55'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +10.58%  '
$ws.Range("E47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5380'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.41%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.842'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.38%  '
$ws.Range("E49").ClearFormats()
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06455'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.97%  '
$ws.Range("E50").ClearFormats()
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.059'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.43%  '
$ws.Range("E51").ClearFormats()
